# Fruta / hortaliza, semanal
# Insert a new weekly record as row 22 (pushing the existing rows 22-33
# down to 23-34) on the Cebollín sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 22; everything below shifts down one row.
$ws.Rows.Item(22).Insert()

# Populate the new row 22 with the new weekly observation.
$ws.Range("A22").Value = 7
$ws.Range("B22").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C22").Value = "Ñuble"
$ws.Range("D22").Value = 44845
$ws.Range("E22").Value = 16
$ws.Range("F22").Value = 100112037
$ws.Range("G22").Value = "Cebollín"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 7500
$ws.Range("L22").Value = 8000
$ws.Range("M22").Value = 7750
$ws.Range("N22").Value = "`$/docena de atados"
$ws.Range("O22").Value = "Provincia de Diguillín"
$ws.Range("P22").Value = 2583
$ws.Range("Q22").Value = 3
$ws.Range("R22").Value = "Hortaliza"
